$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first four data rows (rows 2-5), shifting the remaining
# rows up so the table starts again at row 2. This fixes the off-by-four
# year misalignment in the YoY forecast vector data.
$ws.Rows("2:5").Delete()
